$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename the values used in the "Sexo" column (G) from
# Masculino/Femenino to Hombre/Mujer across the whole sheet.
$usedRange = $ws.UsedRange
[void]$usedRange.Replace("Masculino", "Hombre", 1)
[void]$usedRange.Replace("Femenino", "Mujer", 1)

# Update the current selection on the sheet.
[void]$ws.Range("Q103").Select()

# Reposition the workbook window (best effort).
$win = $excel.ActiveWindow
$win.Left = -120
$win.Top = -120
